# Update the cryptos worksheet with refreshed price/volume data
# as produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.373.11"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.513.04"
$ws.Range("E3").Value = "  +1.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'591.62"
$ws.Range("E5").Value = "  +1.21%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'176.51"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.70%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.512.56"
$ws.Range("E9").Value = "  +1.45%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +3.55%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.94%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.53%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "'0.338"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.969.48"
$ws.Range("E14").Value = "  +1.34%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'25.87"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "68.177.54"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +0.21%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.513.80"
$ws.Range("E18").Value = "  +3.89%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'11.01"
$ws.Range("E19").Value = "  +0.61%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.88%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'351.14"
$ws.Range("E21").Value = "  +0.29%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  +3.64%  "

# Row 23 (swap) -> Litecoin
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.36"
$ws.Range("E23").Value = "  +3.26%  "

# Row 24 (swap) -> Dai
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.03%  "

# Row 25 - NEARProtocol
$ws.Range("D25").Value = "'4.22"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26 - SuiNetwork
$ws.Range("D26").Value = "'1.73"
$ws.Range("E26").Value = "  -4.42%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "'9.25"
$ws.Range("E27").Value = "  +0.94%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.639.13"
$ws.Range("E28").Value = "  +1.19%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'0.989"
$ws.Range("E29").Value = "  -1.11%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0900"
$ws.Range("E30").Value = "  -0.95%  "

# Row 31 - Bittensor
$ws.Range("D31").Value = "'510.85"
$ws.Range("E31").Value = "  +1.59%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'7.83"
$ws.Range("E32").Value = "  +0.97%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +2.06%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +0.98%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  +1.15%  "

# Row 37 - Monero
$ws.Range("D37").Value = "'161.64"
$ws.Range("E37").Value = "  +0.32%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "'18.70"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39 - EthereumClassic
$ws.Range("D39").Value = "'18.39"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -0.44%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +3.99%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - PolygonEcosystemToken
$ws.Range("E43").Value = "  +0.23%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "'4.84"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45 - dogwifhat
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'150.82"
$ws.Range("E46").Value = "  +5.96%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +1.20%  "

# Row 49 (swap) -> Optimism
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.60"
$ws.Range("E49").Value = "  +1.31%  "

# Row 50 (swap) -> Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0741"
$ws.Range("E50").Value = "  +0.15%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "'0.578"
$ws.Range("E51").Value = "  -1.10%  "
